# Refresh GSC export for heatlabs.net-Https.xlsx:
# append the 2025-11-04 data point to the "Chart" sheet.
# (Table sheet keeps its existing Issue/Validation/Pages header - untouched,
# its shared-string indices simply shift because of the new string above.)

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

$newRow = 30

# Write the date as literal text (not an auto-converted date serial) so it
# matches the existing "Date" column cells, which are all stored as shared
# strings. Force text format before assigning, then drop the format
# override again so the cell ends up with the sheet's default style - same
# as every other cell in the column.
$dateCell = $chart.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-11-04"
$dateCell.ClearFormats()

$chart.Cells.Item($newRow, 2).Value = 0.0
$chart.Cells.Item($newRow, 3).Value = 107.0
